$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text so number-like values (e.g. "5.15", "65.00")
# are not auto-converted to numeric by Excel's type inference.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "61.371.65"
$ws.Range("E2").Value = "  -1.21%  "
# Row 3
$ws.Range("D3").Value = "2.438.17"
$ws.Range("E3").Value = "  -0.32%  "
# Row 4
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.14%  "
# Row 5
$ws.Range("D5").Value = "574.94"
$ws.Range("E5").Value = "  -0.83%  "
# Row 6
$ws.Range("D6").Value = "140.86"
$ws.Range("E6").Value = "  -1.42%  "
# Row 7
$ws.Range("E7").Value = "  +0.07%  "
# Row 8
$ws.Range("D8").Value = "0.530"
$ws.Range("E8").Value = "  +0.26%  "
# Row 9
$ws.Range("D9").Value = "2.427.80"
$ws.Range("E9").Value = "  -0.69%  "
# Row 10
$ws.Range("E10").Value = "  +2.02%  "
# Row 11
$ws.Range("E11").Value = "  +1.72%  "
# Row 12
$ws.Range("D12").Value = "5.15"
$ws.Range("E12").Value = "  -0.78%  "
# Row 13
$ws.Range("D13").Value = "0.339"
# Row 14
$ws.Range("D14").Value = "26.10"
$ws.Range("E14").Value = "  -0.92%  "
# Row 15
$ws.Range("D15").Value = "2.891.36"
$ws.Range("E15").Value = "  +0.92%  "
# Row 16
$ws.Range("E16").Value = "  -1.12%  "
# Row 17
$ws.Range("D17").Value = "61.335.83"
$ws.Range("E17").Value = "  -1.28%  "
# Row 18
$ws.Range("D18").Value = "2.424.81"
$ws.Range("E18").Value = "  -0.56%  "
# Row 19
$ws.Range("D19").Value = "10.59"
$ws.Range("E19").Value = "  -3.07%  "
# Row 20
$ws.Range("D20").Value = "7.25"
$ws.Range("E20").Value = "  +2.13%  "
# Row 21
$ws.Range("D21").Value = "324.55"
$ws.Range("E21").Value = "  -1.37%  "
# Row 22
$ws.Range("D22").Value = "4.05"
$ws.Range("E22").Value = "  -1.30%  "
# Row 23
$ws.Range("D23").Value = "6.05"
$ws.Range("E23").Value = "  +1.26%  "
# Row 24
$ws.Range("E24").Value = "  +0.13%  "
# Row 25
$ws.Range("E25").Value = "  -2.11%  "
# Row 26
$ws.Range("D26").Value = "65.00"
$ws.Range("E26").Value = "  -0.97%  "
# Row 27
$ws.Range("D27").Value = "8.98"
$ws.Range("E27").Value = "  -4.26%  "
# Row 28
$ws.Range("D28").Value = "575.00"
$ws.Range("E28").Value = "  -7.46%  "
# Row 29
$ws.Range("D29").Value = "2.570.24"
$ws.Range("E29").Value = "  +0.39%  "
# Row 30
$ws.Range("E30").Value = "  -0.05%  "
# Row 31
$ws.Range("D31").Value = "0.0₃0912"
$ws.Range("E31").Value = "  -3.49%  "
# Row 32
$ws.Range("D32").Value = "7.87"
$ws.Range("E32").Value = "  -1.57%  "
# Row 33
$ws.Range("E33").Value = "  -5.14%  "
# Row 34
$ws.Range("D34").Value = "1.85"
$ws.Range("E34").Value = "  -1.22%  "
# Row 35
$ws.Range("E35").Value = "  -6.66%  "
# Row 36
$ws.Range("E36").Value = "  +0.15%  "
# Row 37
$ws.Range("D37").Value = "4.64"
$ws.Range("E37").Value = "  -5.65%  "
# Row 38
$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D38").Value = "151.80"
$ws.Range("E38").Value = "  +1.17%  "
# Row 39
$ws.Range("B39").Value = "PolygonEcosystemToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D39").Value = "0.370"
$ws.Range("E39").Value = "  -1.36%  "
# Row 40
$ws.Range("D40").Value = "1.38"
$ws.Range("E40").Value = "  -3.42%  "
# Row 41
$ws.Range("D41").Value = "18.30"
$ws.Range("E41").Value = "  +0.00%  "
# Row 42
$ws.Range("D42").Value = "5.13"
$ws.Range("E42").Value = "  -2.13%  "
# Row 43
$ws.Range("E43").Value = "  +0.01%  "
# Row 44
$ws.Range("D44").Value = "41.73"
$ws.Range("E44").Value = "  -2.52%  "
# Row 45
$ws.Range("E45").Value = "  -5.28%  "
# Row 46
$ws.Range("D46").Value = "2.34"
$ws.Range("E46").Value = "  -5.11%  "
# Row 47
$ws.Range("E47").Value = "  +25.15%  "
# Row 48
$ws.Range("D48").Value = "141.77"
$ws.Range("E48").Value = "  -0.77%  "
# Row 49
$ws.Range("D49").Value = "3.53"
# Row 50
$ws.Range("D50").Value = "0.594"
$ws.Range("E50").Value = "  -0.87%  "
# Row 51
$ws.Range("D51").Value = "0.0508"
$ws.Range("E51").Value = "  -2.99%  "

# Restore default "Normal" style on column D so no stray per-cell style refs remain
# (only the NumberFormat needed to be forced to keep the text values as strings).
$priceRange.Style = "Normal"
